# Pushing changes to xgboost output
# Updates the "RESULTS for $" results table on slide 12 (the Table 36
# graphic frame) with refreshed Test RMSE / Test r2 figures.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# "Table 36" is the 5th shape on this slide.
$shape = $s.Shapes.Item(5)
$tbl = $shape.Table

# Row 2 = "Test RMSE"; columns 2-5 are Base Model / After 1st Tuning /
# After 2nd Tuning / Final Tuning.
$rmseRow = $tbl.Rows.Item(2)
$rmseRow.Cells(2).Shape.TextFrame.TextRange.Text = "360.56"
$rmseRow.Cells(3).Shape.TextFrame.TextRange.Text = "295.44"
$rmseRow.Cells(4).Shape.TextFrame.TextRange.Text = "257.61"
$rmseRow.Cells(5).Shape.TextFrame.TextRange.Text = "252.76"

# Row 3 = "Test r2"
$r2Row = $tbl.Rows.Item(3)
$r2Row.Cells(2).Shape.TextFrame.TextRange.Text = "79.07%"
$r2Row.Cells(3).Shape.TextFrame.TextRange.Text = "85.95%"
$r2Row.Cells(4).Shape.TextFrame.TextRange.Text = "89.32%"
$r2Row.Cells(5).Shape.TextFrame.TextRange.Text = "89.67%"
